# Auto-generated: apply updated market-data values to Chocobo Profits sheets
# Mirrors a scheduled-runner refresh of currentAveragePrice / LevePrice / LeveProfit columns
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 533.3333
$ws.Range("J32").Value = 433.33334
$ws.Range("L32").Value = 433.33334
$ws.Range("N32").Value = -1085.33334
$ws.Range("H33").Value = 162.5
$ws.Range("I33").Value = 101.76471
$ws.Range("K33").Value = 101.76471
$ws.Range("M33").Value = 127.23529
$ws.Range("H137").Value = 2954.923
$ws.Range("I137").Value = 1466.5555
$ws.Range("J137").Value = 4230.6665
$ws.Range("K137").Value = 4399.666499999999
$ws.Range("L137").Value = 12691.9995
$ws.Range("M137").Value = -1849.666499999999
$ws.Range("N137").Value = -17791.9995
$ws.Range("H138").Value = 4303.78
$ws.Range("I138").Value = 790.871
$ws.Range("J138").Value = 5882.0435
$ws.Range("K138").Value = 2372.613
$ws.Range("L138").Value = 17646.1305
$ws.Range("M138").Value = 2767.387
$ws.Range("N138").Value = -27926.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 656.85
$ws.Range("I2").Value = 664.25
$ws.Range("J2").Value = 627.25
$ws.Range("K2").Value = 664.25
$ws.Range("L2").Value = 627.25
$ws.Range("M2").Value = -551.25
$ws.Range("N2").Value = -853.25
$ws.Range("H32").Value = 3142.553
$ws.Range("I32").Value = 2799.039
$ws.Range("K32").Value = 2799.039
$ws.Range("M32").Value = -2512.039
$ws.Range("H52").Value = 27390
$ws.Range("J52").Value = 27390
$ws.Range("L52").Value = 27390
$ws.Range("N52").Value = -28026
$ws.Range("H61").Value = 1747.0555
$ws.Range("I61").Value = 1269.4166
$ws.Range("J61").Value = 2702.3333
$ws.Range("K61").Value = 1269.4166
$ws.Range("L61").Value = 2702.3333
$ws.Range("M61").Value = -1057.4166
$ws.Range("N61").Value = -3126.3333
$ws.Range("H109").Value = 25604.762
$ws.Range("J109").Value = 25604.762
$ws.Range("L109").Value = 25604.762
$ws.Range("N109").Value = -28378.762
$ws.Range("H110").Value = 2056.3462
$ws.Range("I110").Value = 1959.9474
$ws.Range("J110").Value = 2318
$ws.Range("K110").Value = 1959.9474
$ws.Range("L110").Value = 2318
$ws.Range("M110").Value = 85.05259999999998
$ws.Range("N110").Value = -6408
$ws.Range("H116").Value = 656.85
$ws.Range("I116").Value = 664.25
$ws.Range("J116").Value = 627.25
$ws.Range("K116").Value = 664.25
$ws.Range("L116").Value = 627.25
$ws.Range("M116").Value = 1629.75
$ws.Range("N116").Value = -5215.25
$ws.Range("H122").Value = 2619.15
$ws.Range("I122").Value = 1745.4615
$ws.Range("J122").Value = 4241.7144
$ws.Range("K122").Value = 5236.3845
$ws.Range("L122").Value = 12725.1432
$ws.Range("M122").Value = -2786.3845
$ws.Range("N122").Value = -17625.1432
$ws.Range("H132").Value = 2074.0508
$ws.Range("I132").Value = 1292.8182
$ws.Range("J132").Value = 4365.6665
$ws.Range("K132").Value = 3878.4546
$ws.Range("L132").Value = 13096.9995
$ws.Range("M132").Value = -1348.4546
$ws.Range("N132").Value = -18156.9995
$ws.Range("H136").Value = 1747.0555
$ws.Range("I136").Value = 1269.4166
$ws.Range("J136").Value = 2702.3333
$ws.Range("K136").Value = 3808.2498
$ws.Range("L136").Value = 8106.999899999999
$ws.Range("M136").Value = -1258.2498
$ws.Range("N136").Value = -13206.9999
$ws.Range("H137").Value = 39786
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 39786
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 39786
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -49986

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 656.85
$ws.Range("I3").Value = 664.25
$ws.Range("J3").Value = 627.25
$ws.Range("K3").Value = 664.25
$ws.Range("L3").Value = 627.25
$ws.Range("M3").Value = -550.25
$ws.Range("N3").Value = -855.25
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50546
$ws.Range("H107").Value = 2556
$ws.Range("I107").Value = 2660
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 2660
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -740
$ws.Range("N107").Value = -6240
$ws.Range("H134").Value = 1732
$ws.Range("I134").Value = 1111.5178
$ws.Range("J134").Value = 3068.423
$ws.Range("K134").Value = 3334.5534
$ws.Range("L134").Value = 9205.269
$ws.Range("M134").Value = -799.5534000000002
$ws.Range("N134").Value = -14275.269
$ws.Range("H140").Value = 58916.25
$ws.Range("J140").Value = 58916.25
$ws.Range("L140").Value = 58916.25
$ws.Range("N140").Value = -69276.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6946791.5
$ws.Range("I31").Value = 1268.7391
$ws.Range("J31").Value = 19235024
$ws.Range("K31").Value = 1268.7391
$ws.Range("L31").Value = 19235024
$ws.Range("M31").Value = -973.7391
$ws.Range("N31").Value = -19235614
$ws.Range("H34").Value = 6946791.5
$ws.Range("I34").Value = 1268.7391
$ws.Range("J34").Value = 19235024
$ws.Range("K34").Value = 1268.7391
$ws.Range("L34").Value = 19235024
$ws.Range("M34").Value = -1066.7391
$ws.Range("N34").Value = -19235428
$ws.Range("H107").Value = 752.4375
$ws.Range("I107").Value = 389
$ws.Range("J107").Value = 1219.7142
$ws.Range("K107").Value = 389
$ws.Range("L107").Value = 1219.7142
$ws.Range("M107").Value = 1531
$ws.Range("N107").Value = -5059.7142
$ws.Range("H122").Value = 2350.611
$ws.Range("I122").Value = 1446.1428
$ws.Range("J122").Value = 2926.182
$ws.Range("K122").Value = 4338.428400000001
$ws.Range("L122").Value = 8778.545999999998
$ws.Range("M122").Value = -1888.428400000001
$ws.Range("N122").Value = -13678.546
$ws.Range("H134").Value = 5597.8335
$ws.Range("I134").Value = 10303.182
$ws.Range("K134").Value = 30909.546
$ws.Range("M134").Value = -28374.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1460.3715
$ws.Range("I5").Value = 367.22223
$ws.Range("J5").Value = 2617.8235
$ws.Range("K5").Value = 1101.66669
$ws.Range("L5").Value = 7853.470499999999
$ws.Range("M5").Value = -989.66669
$ws.Range("N5").Value = -8077.470499999999
$ws.Range("H109").Value = 3889
$ws.Range("I109").Value = 826
$ws.Range("J109").Value = 10015
$ws.Range("K109").Value = 2478
$ws.Range("L109").Value = 30045
$ws.Range("M109").Value = -1438
$ws.Range("N109").Value = -32125
$ws.Range("H113").Value = 542.75757
$ws.Range("I113").Value = 487.52777
$ws.Range("K113").Value = 1462.58331
$ws.Range("M113").Value = 707.41669
$ws.Range("H118").Value = 458.1111
$ws.Range("I118").Value = 458.1111
$ws.Range("K118").Value = 1374.3333
$ws.Range("M118").Value = -131.3333
$ws.Range("H131").Value = 873.40845
$ws.Range("J131").Value = 938.1774
$ws.Range("L131").Value = 2814.5322
$ws.Range("N131").Value = -12894.5322
$ws.Range("H135").Value = 1460.3715
$ws.Range("I135").Value = 367.22223
$ws.Range("J135").Value = 2617.8235
$ws.Range("K135").Value = 3305.00007
$ws.Range("L135").Value = 23560.4115
$ws.Range("M135").Value = -770.0000700000001
$ws.Range("N135").Value = -28630.4115

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 8337.333000000001
$ws.Range("I24").Value = 8006
$ws.Range("K24").Value = 8006
$ws.Range("M24").Value = -7833
$ws.Range("H102").Value = 1970.5
$ws.Range("I102").Value = 1596.8125
$ws.Range("J102").Value = 2397.5715
$ws.Range("K102").Value = 1596.8125
$ws.Range("L102").Value = 2397.5715
$ws.Range("M102").Value = 25.1875
$ws.Range("N102").Value = -5641.5715
$ws.Range("H137").Value = 42786
$ws.Range("J137").Value = 42786
$ws.Range("L137").Value = 42786
$ws.Range("N137").Value = -52986

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5572.1816
$ws.Range("I7").Value = 3413.4285
$ws.Range("J7").Value = 9350
$ws.Range("K7").Value = 3413.4285
$ws.Range("L7").Value = 9350
$ws.Range("M7").Value = -3301.4285
$ws.Range("N7").Value = -9574
$ws.Range("H40").Value = 7604.5
$ws.Range("I40").Value = 5023
$ws.Range("J40").Value = 9716.637000000001
$ws.Range("K40").Value = 5023
$ws.Range("L40").Value = 9716.637000000001
$ws.Range("M40").Value = -4887
$ws.Range("N40").Value = -9988.637000000001
$ws.Range("H46").Value = 2589.2
$ws.Range("I46").Value = 2772.5
$ws.Range("J46").Value = 2467
$ws.Range("K46").Value = 2772.5
$ws.Range("L46").Value = 2467
$ws.Range("M46").Value = -2584.5
$ws.Range("N46").Value = -2843
$ws.Range("H61").Value = 1272.3889
$ws.Range("I61").Value = 1162.6923
$ws.Range("J61").Value = 1557.6
$ws.Range("K61").Value = 1162.6923
$ws.Range("L61").Value = 1557.6
$ws.Range("M61").Value = -960.6922999999999
$ws.Range("N61").Value = -1961.6
$ws.Range("H113").Value = 1272.3889
$ws.Range("I113").Value = 1162.6923
$ws.Range("J113").Value = 1557.6
$ws.Range("K113").Value = 1162.6923
$ws.Range("L113").Value = 1557.6
$ws.Range("M113").Value = 1007.3077
$ws.Range("N113").Value = -5897.6
$ws.Range("H126").Value = 5572.1816
$ws.Range("I126").Value = 3413.4285
$ws.Range("J126").Value = 9350
$ws.Range("K126").Value = 10240.2855
$ws.Range("L126").Value = 28050
$ws.Range("M126").Value = -7770.2855
$ws.Range("N126").Value = -32990
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 756.8570999999999
$ws.Range("I107").Value = 649.5
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1948.5
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -28.5
$ws.Range("N107").Value = -6540
$ws.Range("H111").Value = 39800
$ws.Range("J111").Value = 39800
$ws.Range("L111").Value = 39800
$ws.Range("N111").Value = -47980
$ws.Range("H126").Value = 255121.4
$ws.Range("I126").Value = 987.36
$ws.Range("J126").Value = 628847.9399999999
$ws.Range("K126").Value = 2962.08
$ws.Range("L126").Value = 1886543.82
$ws.Range("M126").Value = -492.0799999999999
$ws.Range("N126").Value = -1891483.82
$ws.Range("H132").Value = 6174341.5
$ws.Range("I132").Value = 690.0606
$ws.Range("J132").Value = 15875794
$ws.Range("K132").Value = 2070.1818
$ws.Range("L132").Value = 47627382
$ws.Range("M132").Value = 459.8181999999997
$ws.Range("N132").Value = -47632442

Write-Output "Updated Chocobo_Profits market data across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"
